$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- L4: new cell "De ZIB invariant verwijderen?" ---
$ws.Range("L4").Value = "De ZIB invariant verwijderen?"

# --- L6: new cell "De valueset onder de loop laten nemen bij ZIB" ---
$ws.Range("L6").Value = "De valueset onder de loop laten nemen bij ZIB"

# --- M3: new cell "Reference goed zetten" ---
$ws.Range("M3").Value = "Reference goed zetten"

# --- L3: extend the bold "gForge trackerID 12798" note with " - ZIB --> ZIB-526" ---
$newL3 = "Allowing a Condition reference: gForge trackerID 12798 - ZIB --> ZIB-526"
$cell = $ws.Range("L3")
$cell.Value = $newL3
$boldStart = $newL3.IndexOf("gForge trackerID 12798") + 1
$boldLen = "gForge trackerID 12798 - ZIB --> ZIB-526".Length
$cell.Characters($boldStart, $boldLen).Font.Bold = $true

# --- view state: scroll so column I is leftmost, zoom to 130%, select J19 ---
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$win.Zoom = 130
$ws.Range("J19").Select()
